# Apply the LiteratureReview.xlsx update:
#  - Add a new reference row (row 26) describing the WGS84 source used for
#    the Gravity.py module's normal gravity formulas.
#  - Select/scroll the sheet roughly to where the new row is, mirroring the
#    author's on-screen state when the workbook was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New bibliography entry (row 26) -------------------------------------
# Column order below matches the order shared strings were introduced in
# the source document (Название/B, Номер/A, Дата обращения/G, Описание/C).
$ws.Range("B26").Value = "World Geodetic System (1984) // ahrs.readthedocs.io"
$ws.Range("A26").Value = "[25]"
$ws.Range("G26").Value = "(16.04.2025)"
$ws.Range("C26").Value = "​Сайт предоставляет подробную информацию о модели WGS84 (World Geodetic System 1984), включая формулы для расчета нормального ускорения свободного падения на поверхности Земли и на высотах над ней."

# Row grew tall (wrapped description text), matching the other reference rows.
$ws.Rows(26).RowHeight = 75

# --- View state -------------------------------------------------------
# Scroll/selection moved further down the sheet to show the newly added row.
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C32").Select() | Out-Null
